$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Large triradiate symetrical rays"
$ws.Range("A3").Value = "Small triradiate symetrical rays"
